# Atualizado por script em 02-01-2024 20:45
# Appends 4 new match rows (87-90) to the Azerbaijan Premier League 2023-2024
# sheet, mirroring the existing table's layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 86

$rows = @(
    @{
        Indice = 86
        pais = "azerbaijan"
        torneio = "premier-league"
        temporada = "2023-2024"
        data_partida = 45283.64583333334
        home = "Sumqayit"
        home_ft_gols = 1
        away = "Sabail"
        away_ft_gols = 0
        home_opening_odds = 1.9
        home_opening_data_hora = "22/12/2023 03:42"
        home_closing_odds = 2.14
        home_closing_data_hora = "23/12/2023 15:21"
        draw_opening_odds = 3.19
        draw_opening_data_hora = "22/12/2023 03:42"
        draw_closing_odds = 3.6
        draw_closing_data_hora = "23/12/2023 15:21"
        away_opening_odds = 3.76
        away_opening_data_hora = "22/12/2023 03:42"
        away_closing_odds = 3.06
        away_closing_data_hora = "23/12/2023 15:21"
        url_partida = "https://www.betexplorer.com/football/azerbaijan/premier-league/sumqayit-fk-sabail/EgB5fHfK/"
    },
    @{
        Indice = 87
        pais = "azerbaijan"
        torneio = "premier-league"
        temporada = "2023-2024"
        data_partida = 45284.375
        home = "Turan"
        home_ft_gols = 3
        away = "Araz"
        away_ft_gols = 1
        home_opening_odds = 2.83
        home_opening_data_hora = "22/12/2023 01:12"
        home_closing_odds = 1.88
        home_closing_data_hora = "24/12/2023 08:40"
        draw_opening_odds = 2.99
        draw_opening_data_hora = "22/12/2023 01:12"
        draw_closing_odds = 3.74
        draw_closing_data_hora = "24/12/2023 08:40"
        away_opening_odds = 2.4
        away_opening_data_hora = "22/12/2023 01:12"
        away_closing_odds = 3.64
        away_closing_data_hora = "24/12/2023 08:40"
        url_partida = "https://www.betexplorer.com/football/azerbaijan/premier-league/turan-araz-pfk/8dKFXgm7/"
    },
    @{
        Indice = 88
        pais = "azerbaijan"
        torneio = "premier-league"
        temporada = "2023-2024"
        data_partida = 45284.54166666666
        home = "Neftci Baku"
        home_ft_gols = 2
        away = "Kapaz"
        away_ft_gols = 0
        home_opening_odds = 1.48
        home_opening_data_hora = "23/12/2023 00:12"
        home_closing_odds = 1.83
        home_closing_data_hora = "24/12/2023 12:56"
        draw_opening_odds = 3.85
        draw_opening_data_hora = "23/12/2023 00:12"
        draw_closing_odds = 3.58
        draw_closing_data_hora = "24/12/2023 12:56"
        away_opening_odds = 5.95
        away_opening_data_hora = "23/12/2023 00:12"
        away_closing_odds = 4.05
        away_closing_data_hora = "24/12/2023 12:56"
        url_partida = "https://www.betexplorer.com/football/azerbaijan/premier-league/neftci-baku-kapaz/z3C1ecuE/"
    },
    @{
        Indice = 89
        pais = "azerbaijan"
        torneio = "premier-league"
        temporada = "2023-2024"
        data_partida = 45284.64583333334
        home = "Zira"
        home_ft_gols = 0
        away = "Qarabag"
        away_ft_gols = 1
        home_opening_odds = 5.4
        home_opening_data_hora = "23/12/2023 00:12"
        home_closing_odds = 8.85
        home_closing_data_hora = "24/12/2023 15:27"
        draw_opening_odds = 3.74
        draw_opening_data_hora = "23/12/2023 00:12"
        draw_closing_odds = 4.37
        draw_closing_data_hora = "24/12/2023 15:27"
        away_opening_odds = 1.53
        away_opening_data_hora = "23/12/2023 00:12"
        away_closing_odds = 1.37
        away_closing_data_hora = "24/12/2023 15:25"
        url_partida = "https://www.betexplorer.com/football/azerbaijan/premier-league/zira-fk-qarabag-agdam/MaOJWD2D/"
    }
)

$targetRow = $lastRow
foreach ($row in $rows) {
    $targetRow = $targetRow + 1

    # Clone the formatting (styles only, not values) of the last existing
    # data row so the new row matches the table's look (bold/bordered
    # index column, date-formatted match-date column, etc.).
    $ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
    $ws.Range("A" + $targetRow + ":V" + $targetRow).PasteSpecial(-4122)

    $ws.Range("A" + $targetRow).Value = $row.Indice
    $ws.Range("B" + $targetRow).Value = $row.pais
    $ws.Range("C" + $targetRow).Value = $row.torneio
    $ws.Range("D" + $targetRow).Value = $row.temporada
    $ws.Range("E" + $targetRow).Value = $row.data_partida
    $ws.Range("F" + $targetRow).Value = $row.home
    $ws.Range("G" + $targetRow).Value = $row.home_ft_gols
    $ws.Range("H" + $targetRow).Value = $row.away
    $ws.Range("I" + $targetRow).Value = $row.away_ft_gols
    $ws.Range("J" + $targetRow).Value = $row.home_opening_odds
    $ws.Range("K" + $targetRow).Value = $row.home_opening_data_hora
    $ws.Range("L" + $targetRow).Value = $row.home_closing_odds
    $ws.Range("M" + $targetRow).Value = $row.home_closing_data_hora
    $ws.Range("N" + $targetRow).Value = $row.draw_opening_odds
    $ws.Range("O" + $targetRow).Value = $row.draw_opening_data_hora
    $ws.Range("P" + $targetRow).Value = $row.draw_closing_odds
    $ws.Range("Q" + $targetRow).Value = $row.draw_closing_data_hora
    $ws.Range("R" + $targetRow).Value = $row.away_opening_odds
    $ws.Range("S" + $targetRow).Value = $row.away_opening_data_hora
    $ws.Range("T" + $targetRow).Value = $row.away_closing_odds
    $ws.Range("U" + $targetRow).Value = $row.away_closing_data_hora
    $ws.Range("V" + $targetRow).Value = $row.url_partida
}

Write-Host ("Added rows {0}-{1}" -f ($lastRow + 1), $targetRow)
